# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 3 for the
# "b8e29229-..." entry, on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-12 20:36:32"
$wsZhCn.Range("H3").Value = "2016-03-12 20:36:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-12 20:36:35"
$wsDeDe.Range("H3").Value = "2016-03-12 20:36:55"
